$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# STEP 1: Fix the dependent/independent variable wording.
# ---------------------------------------------------------------------------

# 1a. "...being the dependent variables and the total score" ->
#     "...being the independent variables and the total score"
#     (insert "in" right before "dependent variables and the total score")
$r = $d.Content
$r.Find.Execute("dependent variables and the total score")
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertBefore("in")

# 1b. "...as the independent variable. And assess..." ->
#     "...as the dependent variable. And assess..."
$r = $d.Content
$r.Find.Execute("as the independent variable. And assess")
$r.Text = "as the dependent variable. And assess"

# 1c. "...factor and the independent variable. Histograms..." ->
#     "...factor and the dependent variable. Histograms..."
$r = $d.Content
$r.Find.Execute("factor and the independent variable. Histograms")
$r.Text = "factor and the dependent variable. Histograms"

# ---------------------------------------------------------------------------
# STEP 2: Move the "_GoBack" bookmark so it sits right before the
# (relocated) "dependent variable." sentence instead of right before
# "data distribution.".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$r = $d.Content
$r.Find.Execute("dependent variable. Histograms")
$bmTarget = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $bmTarget)

# ---------------------------------------------------------------------------
# STEP 3: Re-introduce the run boundaries that the diff shows, using a
# harmless formatting no-op (toggle Bold on/off) to force Word to split a
# run at an exact character offset without altering the visible text or
# leaving any formatting mark behind.
# ---------------------------------------------------------------------------
function Split-At([int]$pos) {
    $rr = $d.Range($pos, $pos + 1)
    $rr.Bold = 1
    $rr.Bold = 0
}

# --- Paragraph: "We have planned to do regression analysis ..." ---

# boundary before "in"
$r = $d.Content
$r.Find.Execute("in")
$r.Find.Execute("being the in")
Split-At($r.End - 2)

# boundary before "dependent variables and the total score"
$r = $d.Content
$r.Find.Execute("being the in")
Split-At($r.End)

# boundary before " obtained"
$r = $d.Content
$r.Find.Execute("the total score")
Split-At($r.End)

# boundary before " by the winning team as the "
$r = $d.Content
$r.Find.Execute(" obtained")
Split-At($r.End)

# boundary before "dependent variable. And assess"
$r = $d.Content
$r.Find.Execute("team as the ")
Split-At($r.End)

# --- Paragraph: "So considering a time period from ..." ---

# boundary before "p between that factor and the "
$r = $d.Content
$r.Find.Execute("more significant is the relationshi")
Split-At($r.End)

# boundary before "dependent variable." (bookmark sits here too)
$r = $d.Content
$r.Find.Execute("p between that factor and the ")
Split-At($r.End)

# boundary before " Histograms and scatter plots..."
$r = $d.Content
$r.Find.Execute("dependent variable.")
Split-At($r.End)

# boundary before "the " (right before "data distribution.")
$r = $d.Content
$r.Find.Execute("for better visualization of ")
Split-At($r.End)

# boundary before "data distribution."
$r = $d.Content
$r.Find.Execute("for better visualization of the ")
Split-At($r.End)
